$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "MyCat" rounded-rectangle label shape on the slide.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "MyCat") {
            $target = $sh
        }
    }
}

if ($target -ne $null) {
    $tr = $target.TextFrame.TextRange

    # Extend "MyCat" -> "MyCat/Sharding" in-place (stays in the same run,
    # keeping its existing run formatting).
    $head = $tr.Characters(1, 5)
    [void]$head.InsertAfter("/Sharding")

    # Append a new "-Proxy" run after it (inherits the same run formatting
    # as the text immediately preceding it).
    $whole = $target.TextFrame.TextRange
    [void]$whole.InsertAfter("-Proxy")
}
